$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formato Cotizador")

# Row 23 "Desarrollo y Diseño inicial": update price and freeze the total as a static value
$ws.Range("G23").Value = 1100
$ws.Range("H23").Value = 1100

# Row 24 "Servicios Adicionales y Mantenimiento": remove the line item entirely
$ws.Range("B24").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = ""

# Leave the selection where the user left it
$ws.Range("H31:H33").Select()
